$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every Price/Volume cell as text, even when the
# digits look like a plain number (e.g. "303.12"). A bare .Value assignment
# would let Excel auto-detect those as numbers, so pre-format the handful of
# cells whose new price string parses as a plain float to Text ("@") first.
$textForcedRows = @(5, 6, 7, 9, 10, 11, 13, 14, 19, 20, 22, 23, 26, 27, 29, 30, 34, 35, 36, 37, 45, 46, 48, 51)
foreach ($r in $textForcedRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.081.05"
$ws.Range("E2").Value = "  +1.97%  "

$ws.Range("D3").Value = "2.310.49"
$ws.Range("E3").Value = "  +1.84%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "303.12"
$ws.Range("E5").Value = "  +0.85%  "

$ws.Range("D6").Value = "101.63"
$ws.Range("E6").Value = "  +5.70%  "

$ws.Range("D7").Value = "0.504"
$ws.Range("E7").Value = "  +1.88%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  +5.78%  "

$ws.Range("D10").Value = "35.88"
$ws.Range("E10").Value = "  +7.99%  "

$ws.Range("D11").Value = "0.0795"
$ws.Range("E11").Value = "  +0.88%  "

$ws.Range("E12").Value = "  +3.62%  "

$ws.Range("D13").Value = "17.96"
$ws.Range("E13").Value = "  +14.55%  "

$ws.Range("D14").Value = "6.92"
$ws.Range("E14").Value = "  +3.83%  "

$ws.Range("D15").Value = "2.687.33"
$ws.Range("E15").Value = "  +2.59%  "

$ws.Range("D16").Value = "2.326.16"
$ws.Range("E16").Value = "  +2.50%  "

$ws.Range("E17").Value = "  +3.86%  "

$ws.Range("D18").Value = "43.013.19"
$ws.Range("E18").Value = "  +2.08%  "

$ws.Range("D19").Value = "12.61"
$ws.Range("E19").Value = "  +7.76%  "

$ws.Range("D20").Value = "6.18"
$ws.Range("E20").Value = "  +3.06%  "

$ws.Range("D21").Value = "0.0₃0906"
$ws.Range("E21").Value = "  +1.84%  "

$ws.Range("D22").Value = "67.90"
$ws.Range("E22").Value = "  +2.34%  "

$ws.Range("D23").Value = "237.32"
$ws.Range("E23").Value = "  +0.98%  "

$ws.Range("E24").Value = "  +12.49%  "

$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("D27").Value = "24.87"
$ws.Range("E27").Value = "  +3.86%  "

$ws.Range("E28").Value = "  +2.57%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "168.11"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "34.54"
$ws.Range("E30").Value = "  +2.48%  "

$ws.Range("E31").Value = "  +0.90%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D34").Value = "4.73"
$ws.Range("E34").Value = "  +3.21%  "

$ws.Range("D35").Value = "17.16"
$ws.Range("E35").Value = "  +3.37%  "

$ws.Range("D36").Value = "2.41"
$ws.Range("E36").Value = "  +3.60%  "

$ws.Range("D37").Value = "0.0693"
$ws.Range("E37").Value = "  +1.12%  "

$ws.Range("E38").Value = "  +3.63%  "

$ws.Range("E39").Value = "  +2.01%  "

$ws.Range("E40").Value = "  +4.29%  "

$ws.Range("E41").Value = "  +1.14%  "

$ws.Range("E42").Value = "  -4.05%  "

$ws.Range("D43").Value = "1.992.01"
$ws.Range("E43").Value = "  +0.96%  "

$ws.Range("E44").Value = "  +4.42%  "

$ws.Range("D45").Value = "10.31"
$ws.Range("E45").Value = "  +8.02%  "

$ws.Range("D46").Value = "17.74"
$ws.Range("E46").Value = "  +1.21%  "

$ws.Range("E47").Value = "  +3.91%  "

$ws.Range("D48").Value = "56.33"
$ws.Range("E48").Value = "  +7.30%  "

$ws.Range("D49").Value = "2.530.22"
$ws.Range("E49").Value = "  +1.48%  "

$ws.Range("E50").Value = "  +3.78%  "

$ws.Range("D51").Value = "4.58"
$ws.Range("E51").Value = "  +2.33%  "
